$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style from an existing header cell (e.g. F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Update existing numeric values with adjusted precision
$ws.Range("B2").Value = 0.07796894984218589
$ws.Range("D2").Value = 0.1911874935925038

# New numeric data cells
$ws.Range("G2").Value = 0.1256850772835605
$ws.Range("H2").Value = 0.99
